$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 57.35848733333334
$ws.Range("H2").Value = 172.075462
$ws.Range("I2").Value = 0.261658309594631
$ws.Range("J2").Value = 0.261658309594631
$ws.Range("M2").Value = 1.766179333333333
$ws.Range("N2").Value = 5.298538
$ws.Range("O2").Value = 0.006601878454633805
$ws.Range("P2").Value = 0.006601878454633805
$ws.Range("Q2").Value = 101.3053749193951
$ws.Range("R2").Value = 911.748374274556
$ws.Range("S2").Value = 0.001727436356588696
$ws.Range("T2").Value = 0.001727436356588696
$ws.Range("G3").Value = 57.35848733333334
$ws.Range("H3").Value = 172.075462
$ws.Range("I3").Value = 0.261658309594631
$ws.Range("J3").Value = 0.261658309594631
$ws.Range("O3").Value = 0.07933168317417663
$ws.Range("P3").Value = 0.07933168317417665
$ws.Range("Q3").Value = 1217.339271265396
$ws.Range("R3").Value = 10956.05344138856
$ws.Range("S3").Value = 0.02075779411665189
$ws.Range("T3").Value = 0.02075779411665189
$ws.Range("G4").Value = 57.35848733333334
$ws.Range("H4").Value = 172.075462
$ws.Range("I4").Value = 0.261658309594631
$ws.Range("J4").Value = 0.261658309594631
$ws.Range("M4").Value = 6.190911333333333
$ws.Range("N4").Value = 18.572734
$ws.Range("O4").Value = 0.02314127641214326
$ws.Range("P4").Value = 0.02314127641214326
$ws.Range("Q4").Value = 355.1013092947898
$ws.Range("R4").Value = 3195.911783653108
$ws.Range("S4").Value = 0.006055107267863512
$ws.Range("T4").Value = 0.006055107267863513
$ws.Range("G5").Value = 57.35848733333334
$ws.Range("H5").Value = 172.075462
$ws.Range("I5").Value = 0.261658309594631
$ws.Range("J5").Value = 0.261658309594631
$ws.Range("M5").Value = 238.3463463333333
$ws.Range("N5").Value = 715.039039
$ws.Range("O5").Value = 0.8909251619590463
$ws.Range("P5").Value = 0.8909251619590463
$ws.Range("Q5").Value = 13671.18588710678
$ws.Range("R5").Value = 123040.672983961
$ws.Range("S5").Value = 0.2331179718535269
$ws.Range("T5").Value = 0.2331179718535269
$ws.Range("I6").Value = 0.2957894889638607
$ws.Range("J6").Value = 0.2957894889638607
$ws.Range("M6").Value = 1.766179333333333
$ws.Range("N6").Value = 5.298538
$ws.Range("O6").Value = 0.006601878454633805
$ws.Range("P6").Value = 0.006601878454633805
$ws.Range("Q6").Value = 114.5198297853525
$ws.Range("R6").Value = 1030.678468068172
$ws.Range("S6").Value = 0.001952766254297656
$ws.Range("T6").Value = 0.001952766254297656
$ws.Range("I7").Value = 0.2957894889638607
$ws.Range("J7").Value = 0.2957894889638607
$ws.Range("O7").Value = 0.07933168317417663
$ws.Range("P7").Value = 0.07933168317417665
$ws.Range("S7").Value = 0.02346547802473261
$ws.Range("T7").Value = 0.02346547802473262
$ws.Range("I8").Value = 0.2957894889638607
$ws.Range("J8").Value = 0.2957894889638607
$ws.Range("M8").Value = 6.190911333333333
$ws.Range("N8").Value = 18.572734
$ws.Range("O8").Value = 0.02314127641214326
$ws.Range("P8").Value = 0.02314127641214326
$ws.Range("Q8").Value = 401.4213611997552
$ws.Range("R8").Value = 3612.792250797796
$ws.Range("S8").Value = 0.006844946323919299
$ws.Range("T8").Value = 0.0068449463239193
$ws.Range("I9").Value = 0.2957894889638607
$ws.Range("J9").Value = 0.2957894889638607
$ws.Range("M9").Value = 238.3463463333333
$ws.Range("N9").Value = 715.039039
$ws.Range("O9").Value = 0.8909251619590463
$ws.Range("P9").Value = 0.8909251619590463
$ws.Range("Q9").Value = 15454.47990297739
$ws.Range("R9").Value = 139090.3191267965
$ws.Range("S9").Value = 0.2635262983609111
$ws.Range("T9").Value = 0.2635262983609111
$ws.Range("G10").Value = 29.294891
$ws.Range("H10").Value = 87.88467299999999
$ws.Range("I10").Value = 0.1336376186888105
$ws.Range("J10").Value = 0.1336376186888105
$ws.Range("M10").Value = 1.766179333333333
$ws.Range("N10").Value = 5.298538
$ws.Range("O10").Value = 0.006601878454633805
$ws.Range("P10").Value = 0.006601878454633805
$ws.Range("Q10").Value = 51.74003105645266
$ws.Range("R10").Value = 465.6602795080739
$ws.Range("S10").Value = 0.0008822593155502259
$ws.Range("T10").Value = 0.0008822593155502259
$ws.Range("G11").Value = 29.294891
$ws.Range("H11").Value = 87.88467299999999
$ws.Range("I11").Value = 0.1336376186888105
$ws.Range("J11").Value = 0.1336376186888105
$ws.Range("O11").Value = 0.07933168317417663
$ws.Range("P11").Value = 0.07933168317417665
$ws.Range("Q11").Value = 621.7357346698136
$ws.Range("R11").Value = 5595.621612028322
$ws.Range("S11").Value = 0.01060169722597214
$ws.Range("T11").Value = 0.01060169722597214
$ws.Range("G12").Value = 29.294891
$ws.Range("H12").Value = 87.88467299999999
$ws.Range("I12").Value = 0.1336376186888105
$ws.Range("J12").Value = 0.1336376186888105
$ws.Range("M12").Value = 6.190911333333333
$ws.Range("N12").Value = 18.572734
$ws.Range("O12").Value = 0.02314127641214326
$ws.Range("P12").Value = 0.02314127641214326
$ws.Range("Q12").Value = 181.3620727006646
$ws.Range("R12").Value = 1632.258654305982
$ws.Range("S12").Value = 0.003092545073138366
$ws.Range("T12").Value = 0.003092545073138366
$ws.Range("G13").Value = 29.294891
$ws.Range("H13").Value = 87.88467299999999
$ws.Range("I13").Value = 0.1336376186888105
$ws.Range("J13").Value = 0.1336376186888105
$ws.Range("M13").Value = 238.3463463333333
$ws.Range("N13").Value = 715.039039
$ws.Range("O13").Value = 0.8909251619590463
$ws.Range("P13").Value = 0.8909251619590463
$ws.Range("Q13").Value = 6982.330236083249
$ws.Range("R13").Value = 62840.97212474924
$ws.Range("S13").Value = 0.1190611170741497
$ws.Range("T13").Value = 0.1190611170741497
$ws.Range("G14").Value = 67.71760166666667
$ws.Range("H14").Value = 203.152805
$ws.Range("I14").Value = 0.3089145827526977
$ws.Range("J14").Value = 0.3089145827526977
$ws.Range("M14").Value = 1.766179333333333
$ws.Range("N14").Value = 5.298538
$ws.Range("O14").Value = 0.006601878454633805
$ws.Range("P14").Value = 0.006601878454633805
$ws.Range("Q14").Value = 119.6014285665656
$ws.Range("R14").Value = 1076.41285709909
$ws.Range("S14").Value = 0.002039416528197227
$ws.Range("T14").Value = 0.002039416528197227
$ws.Range("G15").Value = 67.71760166666667
$ws.Range("H15").Value = 203.152805
$ws.Range("I15").Value = 0.3089145827526977
$ws.Range("J15").Value = 0.3089145827526977
$ws.Range("O15").Value = 0.07933168317417663
$ws.Range("P15").Value = 0.07933168317417665
$ws.Range("Q15").Value = 1437.194383904784
$ws.Range("R15").Value = 12934.74945514306
$ws.Range("S15").Value = 0.02450671380681999
$ws.Range("T15").Value = 0.02450671380681999
$ws.Range("G16").Value = 67.71760166666667
$ws.Range("H16").Value = 203.152805
$ws.Range("I16").Value = 0.3089145827526977
$ws.Range("J16").Value = 0.3089145827526977
$ws.Range("M16").Value = 6.190911333333333
$ws.Range("N16").Value = 18.572734
$ws.Range("O16").Value = 0.02314127641214326
$ws.Range("P16").Value = 0.02314127641214326
$ws.Range("Q16").Value = 419.2336676243189
$ws.Range("R16").Value = 3773.10300861887
$ws.Range("S16").Value = 0.007148677747222081
$ws.Range("T16").Value = 0.007148677747222082
$ws.Range("G17").Value = 67.71760166666667
$ws.Range("H17").Value = 203.152805
$ws.Range("I17").Value = 0.3089145827526977
$ws.Range("J17").Value = 0.3089145827526977
$ws.Range("M17").Value = 238.3463463333333
$ws.Range("N17").Value = 715.039039
$ws.Range("O17").Value = 0.8909251619590463
$ws.Range("P17").Value = 0.8909251619590463
$ws.Range("Q17").Value = 16140.24293970604
$ws.Range("R17").Value = 145262.1864573544
$ws.Range("S17").Value = 0.2752197746704584
$ws.Range("T17").Value = 0.2752197746704584
